# Alarms.xlsx help-doc update:
#  - "Event Button" action now supports overwriting the button value
#    (<button number>[><value>],..,<button number>[><value>])
#  - widen the "Command" column on the Actions sheet so the longer
#    example text fits
#  - the Actions sheet becomes the active / selected sheet (with C6
#    selected), matching the author having just edited that cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actions")

# --- B6: "<button number>" -> "<button number>[><value>],..,<button number>[><value>]"
$newCommand = "<button number>[><value>],..,<button number>[><value>]"
$ws.Range("B6").Value = $newCommand

# Re-apply the (unchanged) default font to the trailing "[><value>]" portion so
# it carries explicit run formatting, matching the source edit.
$suffix = "[><value>]"
$startPos = $newCommand.Length - $suffix.Length + 1
$chars = $ws.Range("B6").Characters($startPos, $suffix.Length)
$chars.Font.Name = "Calibri"
$chars.Font.Size = 11
$chars.Font.Color = 0

# --- C6: describe the new optional-value overwrite syntax
# (uses the typographic single quotes seen in the original edit)
$lq = [char]8216
$rq = [char]8217
$ws.Range("C6").Value = "triggers the button optional overwriting the button value with <value> (eg. " + $lq + "1>10,2,3>100" + $rq + "); the button number comes from the Events Buttons configuration"

# Widen column B (Command) to fit the longer example text.
$ws.Columns.Item(2).ColumnWidth = 36.5

# Make Actions the active sheet/tab and select the edited cell.
$ws.Activate()
$ws.Range("C6").Select()
